$wb = $excel.ActiveWorkbook

# weibull
$ws = $wb.Worksheets.Item("weibull")
$ws.Range("B2").Value = -3.13810935882744
$ws.Range("C2").Value = 0.188580141371055
$ws.Range("B3").Value = 0.150618291702876
$ws.Range("C3").Value = 0.122038959982672

# lognormal
$ws = $wb.Worksheets.Item("lognormal")
$ws.Range("B2").Value = 2.51392827451678
$ws.Range("C2").Value = 0.201082049136534
$ws.Range("B3").Value = -1.04507703041139
$ws.Range("C3").Value = 0.11006605620096

# llogis
$ws = $wb.Worksheets.Item("llogis")
$ws.Range("B2").Value = -2.36474220234867
$ws.Range("C2").Value = 0.137854567634377
$ws.Range("B3").Value = 0.55764737964053
$ws.Range("C3").Value = 0.111678500072081

# gompertz
$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("B2").Value = -2.70193626478097
$ws.Range("C2").Value = 0.156444834207228
$ws.Range("B3").Value = -0.0190879804587354
$ws.Range("C3").Value = 0.019550634465166

# exp - no value changes

# weibull cov
$ws = $wb.Worksheets.Item("weibull cov")
$ws.Range("A2").Value = 0.0355624697195271
$ws.Range("B2").Value = -0.0152701894491901
$ws.Range("A3").Value = -0.0152701894491901
$ws.Range("B3").Value = 0.0148935077536523

# lognormal cov
$ws = $wb.Worksheets.Item("lognormal cov")
$ws.Range("A2").Value = 0.0404339904849474
$ws.Range("B2").Value = -0.0187037178595124
$ws.Range("A3").Value = -0.0187037178595124
$ws.Range("B3").Value = 0.0121145367276329

# llogis cov
$ws = $wb.Worksheets.Item("llogis cov")
$ws.Range("A2").Value = 0.019003881817661
$ws.Range("B2").Value = 0.00905882053789058
$ws.Range("A3").Value = 0.00905882053789058
$ws.Range("B3").Value = 0.0124720873783497

# gompertz cov
$ws = $wb.Worksheets.Item("gompertz cov")
$ws.Range("A2").Value = 0.0244749861501269
$ws.Range("B2").Value = -0.00172192507099017
$ws.Range("A3").Value = -0.00172192507099017
$ws.Range("B3").Value = 0.000382227307990537

# exp cov - no value changes

$wb.Save()
